$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "06.03"
$ws.Range("A14").Style = $ws.Range("A13").Style

$ws.Range("B14").Value = 2394
$ws.Range("C14").Value = 462
$ws.Range("D14").Value = 1060
$ws.Range("E14").Value = 3916
$ws.Range("F14").Value = 523
$ws.Range("G14").Value = 197
$ws.Range("H14").Value = 4636
$ws.Range("I14").Value = 36359
$ws.Range("J14").Value = 4.25
$ws.Range("K14").Value = 11.28
$ws.Range("L14").Value = 9.970000000000001
